# Auto update stock data
# Updates the "Date_1" (A) and "EBITDA" (B) columns for the most-recent-date
# row of each company block. Values are written as literal text (matching
# the source workbook's inlineStr-as-text storage) rather than being
# auto-converted by Excel into dates/numbers: we briefly force the cell to
# Text number-format, assign the literal string, then restore the cell's
# original style so no visible formatting change is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-LiteralText {
    param($Cell, $Text)
    $origStyle = $Cell.Style
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = $origStyle
}

$updates = @(
    @{ Row = 2;  Date = "2025/11/26"; Ebitda = "4.78" },
    @{ Row = 8;  Date = "2025/11/26"; Ebitda = "7.60" },
    @{ Row = 14; Date = "2025/11/26"; Ebitda = "2.80" },
    @{ Row = 20; Date = "2025/11/26"; Ebitda = "12.39" },
    @{ Row = 26; Date = "2025/11/26"; Ebitda = "9.92" },
    @{ Row = 32; Date = "2025/11/26"; Ebitda = "25.99" },
    @{ Row = 38; Date = "2025/11/26"; Ebitda = $null },
    @{ Row = 44; Date = "2025/11/26"; Ebitda = "10.82" },
    @{ Row = 50; Date = "2025/11/26"; Ebitda = "11.56" },
    @{ Row = 56; Date = "2025/11/26"; Ebitda = "34.30" },
    @{ Row = 62; Date = "2025/11/26"; Ebitda = "11.12" },
    @{ Row = 68; Date = "2025/11/26"; Ebitda = "12.03" },
    @{ Row = 74; Date = "2025/11/26"; Ebitda = "15.51" }
)

foreach ($u in $updates) {
    $dateCell = $ws.Cells.Item($u.Row, 1)
    Set-LiteralText $dateCell $u.Date

    if ($u.Ebitda -ne $null) {
        $ebitdaCell = $ws.Cells.Item($u.Row, 2)
        Set-LiteralText $ebitdaCell $u.Ebitda
    }
}

Write-Host "Updated $($updates.Count) rows"
